# #327 Ajout des profils d'acces
# - Bump the "Date" metadata value on the Metadata sheet.
# - Re-order the two right-most "Mapping" columns on the Elements sheet
#   (the business-mapping column now comes before the RIM mapping column),
#   swapping both the header/data content and the column widths.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Metadata sheet: refresh the generation Date value.
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# ------------------------------------------------------------------
# 2. Elements sheet: swap columns AK (37) and AL (38).
#    Column AK held "Mapping: RIM Mapping", AL held
#    "Mapping: Spécification métier vers l'extension ROR
#    TelecomCommunicationChannel" - they swap places.
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

$colLeft = 37  # AK
$colRight = 38 # AL

$rowCount = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $leftCell = $ws.Cells.Item($r, $colLeft)
    $rightCell = $ws.Cells.Item($r, $colRight)

    $leftVal = $leftCell.Value()
    $rightVal = $rightCell.Value()

    # Skip rows where both sides already hold the same (e.g. both blank)
    # value - nothing to swap, and leaving them alone avoids needlessly
    # rewriting cells that do not actually change.
    if ($leftVal -ne $rightVal) {
        $leftCell.Value = $rightVal
        $rightCell.Value = $leftVal
    }
}

# Swap the column widths to match the new (wider/narrower) contents.
# (target widths are 88.08984375 / 24.98046875 "characters"; the inputs
# below are chosen so the engine's internal pixel-rounding lands as close
# as possible to those exact figures.)
$ws.Range("AK1").ColumnWidth = 87.33333333333333
$ws.Range("AL1").ColumnWidth = 24.166666666666668
